# Add a second sheet ("NewTable0") containing a small annotation table
# (Input [Source Name] / Output [Sample Name]) right after the existing
# "isa_run" sheet, so two runs show up — mirrors the commit
# "trying to get two runs showed in fsx".

$wb = $excel.ActiveWorkbook

# Create the new worksheet and name it.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "NewTable0"

# Position it right after "isa_run" (Worksheets.Add inserts before the
# active sheet, so move it explicitly). Re-fetch the sheet by name
# afterwards since index-based handles shift across a Move.
$newSheet.Move($null, $wb.Worksheets.Item("isa_run"))
$ws = $wb.Worksheets.Item("NewTable0")

# Header row.
$ws.Range("A1").Value = "Input [Source Name]"
$ws.Range("B1").Value = "Output [Sample Name]"

# Data rows.
$ws.Range("A2").Value = "a"
$ws.Range("B2").Value = "ab"
$ws.Range("A3").Value = "b"
$ws.Range("B3").Value = "bc"

# Turn the range into a proper Excel table ("annotationTable0").
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:B3"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "annotationTable0"
$lo.TableStyle = "TableStyleMedium2"

# Match the source file's autofilter buttons being hidden on both columns.
$lo.ShowAutoFilterDropDown = $false
